# Auto-generated edit script: updates cryptos list values (prices/volumes)
# and fixes two row orderings (Chainlink/WrappedEther swap at rows 15-16,
# ImmutableX/Dai swap at rows 25-26), matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.548.67"
$ws.Range("E2").Value = "  +1.93%  "

# Row 3
$ws.Range("D3").Value = "2.492.60"
$ws.Range("E3").Value = "  +2.52%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.994"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.59%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.97%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.547"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.995"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.72%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.26%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.52%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.31%  "

# Row 12
$ws.Range("E12").Value = "  +2.29%  "

# Row 13
$ws.Range("D13").Value = "2.875.62"
$ws.Range("E13").Value = "  +2.28%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.90%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.473.76"
$ws.Range("E16").Value = "  +1.60%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.761"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.28%  "

# Row 18
$ws.Range("D18").Value = "41.752.75"
$ws.Range("E18").Value = "  +2.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.88%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  +2.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "70.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.90%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.63"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.92%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.01%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.06%  "

# Row 28
$ws.Range("E28").Value = "  +1.47%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.87%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.68%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "154.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.52%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.93%  "

# Row 34
$ws.Range("E34").Value = "  -1.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0757"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.75%  "

# Row 36
$ws.Range("E36").Value = "  -1.39%  "

# Row 37
$ws.Range("E37").Value = "  -0.68%  "

# Row 38
$ws.Range("E38").Value = "  -1.35%  "

# Row 39
$ws.Range("E39").Value = "  +0.53%  "

# Row 40
$ws.Range("E40").Value = "  +0.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.79"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.31%  "

# Row 44
$ws.Range("D44").Value = "1.948.91"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45
$ws.Range("E45").Value = "  +1.48%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.43%  "

# Row 48
$ws.Range("D48").Value = "2.735.54"
$ws.Range("E48").Value = "  +2.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "96.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

# Row 51
$ws.Range("E51").Value = "  -1.17%  "

